$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (item 7): update dates B12:E12 to 2015-04-02 (serial 42096) ---
$ws.Range("B12").Value = 42096
$ws.Range("C12").Value = 42096
$ws.Range("D12").Value = 42096
$ws.Range("E12").Value = 42096
$ws.Range("B12").Copy()
$ws.Range("C12:E12").PasteSpecial(-4122)
$ws.Range("F12").Value = "Con segunda revisión editor, revisión de la coordinadora, falta confirmación editora gráfica."
$ws.Rows.Item(12).RowHeight = 75.75

# --- Row 13 (item 8): new dates B13:E13 = 2015-04-07 (serial 42101) ---
$ws.Range("B13").Value = 42101
$ws.Range("C13").Value = 42101
$ws.Range("D13").Value = 42101
$ws.Range("E13").Value = 42101
$ws.Range("B13").Copy()
$ws.Range("C13:E13").PasteSpecial(-4122)
$ws.Range("F13").Value = "En el Git la versión original del autor. En marcha primera revisión del editor."
$ws.Rows.Item(13).RowHeight = 60.75

# --- New row 25: D25 = "," ---
$ws.Range("D25").Value = ","

# --- Update view: scroll & selection ---
$ws.Range("D25").Select()
